$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2023-08-25 Friday" "2023-08-26 Saturday"

Replace-Text "19×65=1235" "27×49=1323"
Replace-Text "12×15=180" "25×63=1575"
Replace-Text "96×99=9504" "44×76=3344"
Replace-Text "85×97=8245" "55×84=4620"
Replace-Text "49×94=4606" "32×75=2400"

Replace-Text "84×86=7224" "61×27=1647"
Replace-Text "88×91=8008" "75×26=1950"
Replace-Text "58×71=4118" "60×31=1860"
Replace-Text "19×87=1653" "31×95=2945"
Replace-Text "83×88=7304" "27×61=1647"

Replace-Text "93×49=4557" "91×40=3640"
Replace-Text "77×36=2772" "41×68=2788"
Replace-Text "34×40=1360" "12×35=420"
Replace-Text "93×15=1395" "46×64=2944"
Replace-Text "13×71=923" "54×31=1674"

Replace-Text "76×53=4028" "55×12=660"
Replace-Text "54×87=4698" "60×75=4500"
Replace-Text "16×46=736" "58×53=3074"
Replace-Text "91×50=4550" "19×66=1254"
Replace-Text "90×56=5040" "36×15=540"

Replace-Text "76×63=4788" "63×56=3528"
Replace-Text "53×85=4505" "28×28=784"
Replace-Text "30×98=2940" "49×91=4459"
Replace-Text "11×56=616" "50×92=4600"
Replace-Text "95×66=6270" "51×74=3774"
